# Update cryptocurrency price/volume table with the latest scraped values
# (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "315.72", "43.197.43").
# Excel auto-converts plain numeric-looking assignments to the Number type
# (dropping significant trailing zeros, e.g. "0.860" -> 0.86), but the source
# data is text. Force text format, assign, then restore the "Normal" style so
# no stray number-format style is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.197.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.83%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.231.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.78%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.80%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  -7.95%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -9.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0823"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.24%  "

$ws.Range("E14").Value = "  -1.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.572.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.860"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.240.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.136.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0960"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.68%  "

$ws.Range("E22").Value = "  -10.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.54"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -9.69%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "238.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -10.39%  "

$ws.Range("E27").Value = "  -0.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.57%  "

$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0871"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -10.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "154.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.21%  "

$ws.Range("E36").Value = "  -5.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("E39").Value = "  -6.55%  "

$ws.Range("E40").Value = "  -5.39%  "

$ws.Range("E41").Value = "  -11.46%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0323"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.96%  "

$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.799.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.205"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "77.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.68%  "

$ws.Range("E50").Value = "  -8.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "59.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -15.58%  "
